$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 29
$ws1.Range("F6").Value = 553
$ws1.Range("F7").Value = 1719
$ws1.Range("F8").Value = 35
$ws1.Range("F10").Value = 37
$ws1.Range("F11").Value = 1801
$ws1.Range("F12").Value = 132
$ws1.Range("F13").Value = 114
$ws1.Range("F14").Value = 425
$ws1.Range("F17").Value = 2
$ws1.Range("F18").Value = 16
$ws1.Range("F19").Value = 35
$ws1.Range("F21").Value = 52
$ws1.Range("F22").Value = 767
$ws1.Range("F23").Value = 311
$ws1.Range("F24").Value = 169
$ws1.Range("F25").Value = 246
$ws1.Range("F26").Value = 261

# Sheet "全部类型" (sheet4): row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 29
$ws4.Range("F6").Value = 553
$ws4.Range("F7").Value = 1719
$ws4.Range("F9").Value = 35
$ws4.Range("F11").Value = 37
$ws4.Range("F12").Value = 1801
$ws4.Range("F13").Value = 132
$ws4.Range("F14").Value = 114
$ws4.Range("F15").Value = 425
$ws4.Range("F18").Value = 2
$ws4.Range("F19").Value = 16
$ws4.Range("F20").Value = 35
$ws4.Range("F22").Value = 52
$ws4.Range("F23").Value = 767
$ws4.Range("F24").Value = 311
$ws4.Range("F25").Value = 169
$ws4.Range("F26").Value = 246
$ws4.Range("F27").Value = 261
